$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly data record needs to be inserted as row 73, pushing the
# existing rows 73-108 down to 74-109 (dimension grows from R108 to R109).
# Copy row 72 first so the new row inherits the same cell formatting/style
# (e.g. the date style applied to column D) as the rest of the table.
$ws.Rows.Item(72).Copy()
$ws.Rows.Item(73).Insert()

# Populate the new row 73 with the new record's values.
$ws.Cells.Item(73, 1).Value = 4
$ws.Cells.Item(73, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(73, 3).Value = "Los Lagos"
$ws.Cells.Item(73, 4).Value = 45086
$ws.Cells.Item(73, 5).Value = 10
$ws.Cells.Item(73, 6).Value = 100112043
$ws.Cells.Item(73, 7).Value = "Pepino dulce"
$ws.Cells.Item(73, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(73, 9).Value = "Primera"
$ws.Cells.Item(73, 10).Value = 60
$ws.Cells.Item(73, 11).Value = 17000
$ws.Cells.Item(73, 12).Value = 17000
$ws.Cells.Item(73, 13).Value = 17000
$ws.Cells.Item(73, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(73, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(73, 16).Value = 944
$ws.Cells.Item(73, 17).Value = 18
$ws.Cells.Item(73, 18).Value = "Hortaliza"
